$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 11.915632
$ws.Range("H2").Value = 35.746896
$ws.Range("I2").Value = 0.2203762099850903
$ws.Range("J2").Value = 0.2203762099850904
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 1394.659940308816
$ws.Range("R2").Value = 12551.93946277934
$ws.Range("S2").Value = 0.07152013809551983
$ws.Range("T2").Value = 0.07152013809551984
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 11.915632
$ws.Range("H3").Value = 35.746896
$ws.Range("I3").Value = 0.2203762099850903
$ws.Range("J3").Value = 0.2203762099850904
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 1210.390343410261
$ws.Range("R3").Value = 10893.51309069235
$ws.Range("S3").Value = 0.06207053204024574
$ws.Range("T3").Value = 0.06207053204024573
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 11.915632
$ws.Range("H4").Value = 35.746896
$ws.Range("I4").Value = 0.2203762099850903
$ws.Range("J4").Value = 0.2203762099850904
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 1692.338955837525
$ws.Range("R4").Value = 15231.05060253773
$ws.Range("S4").Value = 0.08678553984932477
$ws.Range("T4").Value = 0.08678553984932477
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 29.800487
$ws.Range("H5").Value = 89.40146100000001
$ws.Range("I5").Value = 0.5511514941691683
$ws.Range("J5").Value = 0.5511514941691684
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 3487.984978102181
$ws.Range("R5").Value = 31391.86480291963
$ws.Range("S5").Value = 0.1788688124602827
$ws.Range("T5").Value = 0.1788688124602828
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 29.800487
$ws.Range("H6").Value = 89.40146100000001
$ws.Range("I6").Value = 0.5511514941691683
$ws.Range("J6").Value = 0.5511514941691684
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 3027.134582011515
$ws.Range("R6").Value = 27244.21123810363
$ws.Range("S6").Value = 0.155235751082983
$ws.Range("T6").Value = 0.155235751082983
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 29.800487
$ws.Range("H7").Value = 89.40146100000001
$ws.Range("I7").Value = 0.5511514941691683
$ws.Range("J7").Value = 0.5511514941691684
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 4232.467489179739
$ws.Range("R7").Value = 38092.20740261766
$ws.Range("S7").Value = 0.2170469306259026
$ws.Range("T7").Value = 0.2170469306259026
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 12.35338333333333
$ws.Range("H8").Value = 37.06015
$ws.Range("I8").Value = 0.2284722958457413
$ws.Range("J8").Value = 0.2284722958457413
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 1445.896353821483
$ws.Range("R8").Value = 13013.06718439335
$ws.Range("S8").Value = 0.07414761398697887
$ws.Range("T8").Value = 0.07414761398697889
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 12.35338333333333
$ws.Range("H9").Value = 37.06015
$ws.Range("I9").Value = 0.2284722958457413
$ws.Range("J9").Value = 0.2284722958457413
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 1254.857140192978
$ws.Range("R9").Value = 11293.7142617368
$ws.Range("S9").Value = 0.06435085239264726
$ws.Range("T9").Value = 0.06435085239264726
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 12.35338333333333
$ws.Range("H10").Value = 37.06015
$ws.Range("I10").Value = 0.2284722958457413
$ws.Range("J10").Value = 0.2284722958457413
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 1754.511372237245
$ws.Range("R10").Value = 15790.6023501352
$ws.Range("S10").Value = 0.08997382946611514
$ws.Range("T10").Value = 0.08997382946611515
